# edit.ps1 - apply the FinalProjectWriteup.docx revision
#
# Summary of the change:
#  1. The sentence "...uses a linked list..." becomes "...uses a doubly
#     linked list..." (the word "doubly " is inserted, and the surrounding
#     text ends up split into three separate runs).
#  2. A large new passage describing the insert/search/delete functions is
#     appended after "...searching an event by keyword. ", and right before
#     "To run this program" the paragraph is split in two (the second
#     paragraph gets a first-line indent of 720 twips / 36pt).
#  3. The two runs " Then in the Final Project..." and "So that is an
#     enhancement option. " are moved so that they come BEFORE the
#     _GoBack bookmark instead of after it.
#
# Note: this engine recomputes/re-coalesces a paragraph's run list whenever
# text is inserted/deleted anywhere inside that paragraph, which would
# silently undo an earlier same-paragraph run split. To avoid that, every
# InsertBefore/InsertParagraphAfter/Delete call for a given paragraph is
# performed first, and the run-splitting (done by toggling Bold on then off
# over the exact sub-range, which forces a run boundary without changing
# any visible formatting) happens only once no more text edits are pending
# for that paragraph.

$d = $word.ActiveDocument

function Split-RunAt($range) {
    $range.Bold = $true
    $range.Bold = $false
}

# ---------------------------------------------------------------------
# Step 1: insert "doubly " before "linked list for each of the dates"
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("uses a linked list for each of the dates", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $find1.Found) {
    throw "Could not find the 'uses a linked list for each of the dates' anchor text"
}
# "uses a " is 7 characters; insert right after it, before "linked"
$doublyStart = $find1.Start + 7
$insPoint1 = $d.Range($doublyStart, $doublyStart)
$insPoint1.InsertBefore("doubly ")
$doublyEnd = $doublyStart + "doubly ".Length

# ---------------------------------------------------------------------
# Step 2: add the new sentences after "...searching an event by keyword. "
#          and split the paragraph before "To run this program"
# ---------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("searching an event by keyword. To run this program", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $find2.Found) {
    throw "Could not find the 'searching an event by keyword. To run this program' anchor text"
}
# "searching an event by keyword. " is 32 characters
$insertAt = $find2.Start + "searching an event by keyword. ".Length

$textD = "The insert function uses a doubly linked list inserting each date at every index which is representative of a month. Instead of using a hash function, each index is assigned a month. The search by date function allows the user to input any date and see if they have any events that day. The search by keyword function on the other asks the user for a keyword and searches through the database for an event with the matching keyword. The delete an event feature allows the user to delete any event, and the delete all allows them to clear the entire calendar. They also have the option to print the contents of the calendar"
$textE = " this is done by traversing through the hash Table"
$textF = ". "

$insPoint2 = $d.Range($insertAt, $insertAt)
$insPoint2.InsertBefore($textD + $textE + $textF)

$dStart = $insertAt
$dEnd = $dStart + $textD.Length
$eEnd = $dEnd + $textE.Length
$fEnd = $eEnd + $textF.Length

# Break the paragraph right before "To run this program" (immediately
# after the newly inserted ". ")
$breakPoint = $d.Range($fEnd, $fEnd)
$breakPoint.InsertParagraphAfter()

# ---------------------------------------------------------------------
# Step 3: move " Then in the Final Project..." and "So that is an
#          enhancement option. " before the _GoBack bookmark
# ---------------------------------------------------------------------
$text1 = " Then in the Final Project when you run it a Main menu should open where all the options are listed. Depending on the choice number you can input a month, date, and event. The only known bugs in this program is the Search by keyword doesn" + [char]8217 + "t work properly. Right now a limitation that this program has is that you can only add one event per date. "
$text2 = "So that is an enhancement option. "

$find4 = $d.Content
$find4.Find.Execute($text1 + $text2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $find4.Found) {
    throw "Could not find the movable 'Then in the Final Project...enhancement option.' text"
}
# Delete the text from its current location (after the bookmark)...
$find4.Delete()

# ...and re-insert it immediately before the _GoBack bookmark.
$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range
$moveInsertStart = $bmRange.Start
$bmInsertPoint = $d.Range($moveInsertStart, $moveInsertStart)
$bmInsertPoint.InsertBefore($text1 + $text2)

$move1End = $moveInsertStart + $text1.Length
$move2End = $move1End + $text2.Length

# ---------------------------------------------------------------------
# All text insertions/deletions are done. Now force the run splits that
# the target markup expects, and set the new paragraph's indent.
# ---------------------------------------------------------------------
Split-RunAt($d.Range($doublyStart, $doublyEnd))
Split-RunAt($d.Range($dStart, $dEnd))
Split-RunAt($d.Range($dEnd, $eEnd))
Split-RunAt($d.Range($eEnd, $fEnd))
Split-RunAt($d.Range($moveInsertStart, $move1End))
Split-RunAt($d.Range($move1End, $move2End))

# Give the new paragraph (the one that now starts with "To run this
# program") a first-line indent of 720 twips (36pt).
$find3 = $d.Content
$find3.Find.Execute("To run this program, first open the main", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $find3.Found) {
    throw "Could not find the 'To run this program, first open the main' text after the paragraph split"
}
$newPara = $find3.Paragraphs.Item(1)
$newPara.Range.ParagraphFormat.FirstLineIndent = 36

Write-Host "Edit complete"
